# Updated cryptos list with latest coinranking.com price/volume snapshot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell is forced to Text ("@") before the write so numeric-looking
# strings (e.g. "7.200", "158.20") keep their exact characters/trailing zeros
# instead of being auto-coerced to a Number by Excel; ClearFormats() afterwards
# drops the temporary number-format override so the cell style is left untouched.
function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "28.002.07"
Set-TextValue "E2" "  +0.21%  "
Set-TextValue "D3" "1.860.33"
Set-TextValue "E3" "  -0.51%  "
Set-TextValue "D4" "1.004"
Set-TextValue "E4" "  +0.31%  "
Set-TextValue "D5" "311.99"
Set-TextValue "E5" "  -0.20%  "
Set-TextValue "D6" "1.003"
Set-TextValue "E6" "  +0.22%  "
Set-TextValue "D7" "0.5085"
Set-TextValue "E7" "  +1.00%  "
Set-TextValue "D8" "0.3818"
Set-TextValue "E8" "  -0.31%  "
Set-TextValue "D9" "0.08298"
Set-TextValue "E9" "  -7.02%  "
Set-TextValue "E10" "  -0.34%  "
Set-TextValue "D11" "41.62"
Set-TextValue "E11" "  +0.17%  "
Set-TextValue "D12" "6.208"
Set-TextValue "E12" "  -2.80%  "
Set-TextValue "D13" "20.54"
Set-TextValue "E13" "  -0.61%  "
Set-TextValue "D14" "1.853.93"
Set-TextValue "E14" "  -0.66%  "
Set-TextValue "D15" "7.200"
Set-TextValue "E15" "  -0.55%  "
Set-TextValue "E16" "  +0.27%  "
Set-TextValue "D17" "0.00001096"
Set-TextValue "E17" "  -0.26%  "
Set-TextValue "D18" "90.66"
Set-TextValue "E18" "  -0.52%  "
Set-TextValue "D19" "0.06628"
Set-TextValue "E19" "  -0.47%  "
Set-TextValue "D20" "17.68"
Set-TextValue "E20" "  -2.47%  "
Set-TextValue "E21" "  +0.23%  "
Set-TextValue "E22" "  -1.52%  "
Set-TextValue "D23" "28.013.01"
Set-TextValue "E23" "  +0.18%  "
Set-TextValue "E24" "  -3.78%  "
Set-TextValue "D25" "2.237"
Set-TextValue "E25" "  -1.88%  "
Set-TextValue "D26" "2.541"
Set-TextValue "E26" "  +1.90%  "
Set-TextValue "D27" "2.074.42"
Set-TextValue "E27" "  -0.07%  "
Set-TextValue "D28" "158.20"
Set-TextValue "E28" "  -0.12%  "
Set-TextValue "D29" "20.46"
Set-TextValue "E29" "  -1.00%  "
Set-TextValue "D30" "124.48"
Set-TextValue "E30" "  -1.39%  "
Set-TextValue "D31" "0.1052"
Set-TextValue "E31" "  -0.89%  "
Set-TextValue "E32" "  -1.84%  "
Set-TextValue "D33" "5.747"
Set-TextValue "E33" "  +2.54%  "
Set-TextValue "D34" "3.590"
Set-TextValue "E34" "  -0.49%  "
Set-TextValue "D35" "9.425"
Set-TextValue "E35" "  -0.94%  "
Set-TextValue "D36" "0.06513"
Set-TextValue "E36" "  -0.63%  "
Set-TextValue "D37" "0.02411"
Set-TextValue "E37" "  +0.49%  "
Set-TextValue "D38" "0.2165"
Set-TextValue "E38" "  -0.74%  "
Set-TextValue "D39" "1.205"
Set-TextValue "E39" "  -0.08%  "
Set-TextValue "D40" "0.6438"
Set-TextValue "E40" "  +1.04%  "
Set-TextValue "D41" "1.221"
Set-TextValue "E41" "  -4.86%  "
Set-TextValue "D42" "4.904"
Set-TextValue "E42" "  -0.05%  "
Set-TextValue "D43" "11.20"
Set-TextValue "E43" "  -2.63%  "
Set-TextValue "D44" "0.6079"
Set-TextValue "E44" "  +1.26%  "
Set-TextValue "D45" "13.10"
Set-TextValue "E45" "  -0.30%  "
Set-TextValue "D46" "1.283"
Set-TextValue "E46" "  +0.31%  "
Set-TextValue "D47" "3.658"
Set-TextValue "E47" "  -0.22%  "
Set-TextValue "D48" "2.013"
Set-TextValue "E48" "  +1.10%  "
Set-TextValue "D49" "1.207"
Set-TextValue "D50" "120.41"
Set-TextValue "E50" "  -0.30%  "
Set-TextValue "D51" "78.42"
Set-TextValue "E51" "  -0.75%  "
